$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 31.824752
$ws.Cells.Item(2, 8).Value = 95.47425600000001
$ws.Cells.Item(2, 9).Value = 0.886907633630525
$ws.Cells.Item(2, 10).Value = 0.886907633630525
$ws.Cells.Item(2, 13).Value = 5.740110333333334
$ws.Cells.Item(2, 14).Value = 17.220331
$ws.Cells.Item(2, 15).Value = 0.2861925343043439
$ws.Cells.Item(2, 16).Value = 0.2861925343043439
$ws.Cells.Item(2, 17).Value = 182.6775878109707
$ws.Cells.Item(2, 18).Value = 1644.098290298736
$ws.Cells.Item(2, 19).Value = 0.2538263433625885
$ws.Cells.Item(2, 20).Value = 0.2538263433625885
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 31.824752
$ws.Cells.Item(3, 8).Value = 95.47425600000001
$ws.Cells.Item(3, 9).Value = 0.886907633630525
$ws.Cells.Item(3, 10).Value = 0.886907633630525
$ws.Cells.Item(3, 15).Value = 0.2917347240316885
$ws.Cells.Item(3, 16).Value = 0.2917347240316885
$ws.Cells.Item(3, 17).Value = 186.215184810288
$ws.Cells.Item(3, 18).Value = 1675.936663292592
$ws.Cells.Item(3, 19).Value = 0.2587417537387992
$ws.Cells.Item(3, 20).Value = 0.2587417537387992
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 31.824752
$ws.Cells.Item(4, 8).Value = 95.47425600000001
$ws.Cells.Item(4, 9).Value = 0.886907633630525
$ws.Cells.Item(4, 10).Value = 0.886907633630525
$ws.Cells.Item(4, 13).Value = 6.759986
$ws.Cells.Item(4, 14).Value = 20.279958
$ws.Cells.Item(4, 15).Value = 0.3370418707750538
$ws.Cells.Item(4, 16).Value = 0.3370418707750538
$ws.Cells.Item(4, 17).Value = 215.134877973472
$ws.Cells.Item(4, 18).Value = 1936.213901761248
$ws.Cells.Item(4, 19).Value = 0.2989250080435082
$ws.Cells.Item(4, 20).Value = 0.2989250080435082
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 31.824752
$ws.Cells.Item(5, 8).Value = 95.47425600000001
$ws.Cells.Item(5, 9).Value = 0.886907633630525
$ws.Cells.Item(5, 10).Value = 0.886907633630525
$ws.Cells.Item(5, 13).Value = 1.705448333333333
$ws.Cells.Item(5, 14).Value = 5.116345
$ws.Cells.Item(5, 15).Value = 0.0850308708889137
$ws.Cells.Item(5, 16).Value = 0.0850308708889137
$ws.Cells.Item(5, 17).Value = 54.27547025714667
$ws.Cells.Item(5, 18).Value = 488.47923231432
$ws.Cells.Item(5, 19).Value = 0.07541452848562916
$ws.Cells.Item(5, 20).Value = 0.07541452848562916
$ws.Cells.Item(6, 9).Value = 0.06502043684278042
$ws.Cells.Item(6, 10).Value = 0.06502043684278042
$ws.Cells.Item(6, 13).Value = 5.740110333333334
$ws.Cells.Item(6, 14).Value = 17.220331
$ws.Cells.Item(6, 15).Value = 0.2861925343043439
$ws.Cells.Item(6, 16).Value = 0.2861925343043439
$ws.Cells.Item(6, 17).Value = 13.39234900057567
$ws.Cells.Item(6, 18).Value = 120.531141005181
$ws.Cells.Item(6, 19).Value = 0.01860836360161086
$ws.Cells.Item(6, 20).Value = 0.01860836360161086
$ws.Cells.Item(7, 9).Value = 0.06502043684278042
$ws.Cells.Item(7, 10).Value = 0.06502043684278042
$ws.Cells.Item(7, 15).Value = 0.2917347240316885
$ws.Cells.Item(7, 16).Value = 0.2917347240316885
$ws.Cells.Item(7, 19).Value = 0.01896871919874838
$ws.Cells.Item(7, 20).Value = 0.01896871919874838
$ws.Cells.Item(8, 9).Value = 0.06502043684278042
$ws.Cells.Item(8, 10).Value = 0.06502043684278042
$ws.Cells.Item(8, 13).Value = 6.759986
$ws.Cells.Item(8, 14).Value = 20.279958
$ws.Cells.Item(8, 15).Value = 0.3370418707750538
$ws.Cells.Item(8, 16).Value = 0.3370418707750538
$ws.Cells.Item(8, 17).Value = 15.771838256362
$ws.Cells.Item(8, 18).Value = 141.946544307258
$ws.Cells.Item(8, 19).Value = 0.02191460967210195
$ws.Cells.Item(8, 20).Value = 0.02191460967210195
$ws.Cells.Item(9, 9).Value = 0.06502043684278042
$ws.Cells.Item(9, 10).Value = 0.06502043684278042
$ws.Cells.Item(9, 13).Value = 1.705448333333333
$ws.Cells.Item(9, 14).Value = 5.116345
$ws.Cells.Item(9, 15).Value = 0.0850308708889137
$ws.Cells.Item(9, 16).Value = 0.0850308708889137
$ws.Cells.Item(9, 17).Value = 3.979010499121667
$ws.Cells.Item(9, 18).Value = 35.811094492095
$ws.Cells.Item(9, 19).Value = 0.005528744370319229
$ws.Cells.Item(9, 20).Value = 0.005528744370319229
$ws.Cells.Item(10, 7).Value = 1.696588
$ws.Cells.Item(10, 8).Value = 5.089764000000001
$ws.Cells.Item(10, 9).Value = 0.04728133775640876
$ws.Cells.Item(10, 10).Value = 0.04728133775640876
$ws.Cells.Item(10, 13).Value = 5.740110333333334
$ws.Cells.Item(10, 14).Value = 17.220331
$ws.Cells.Item(10, 15).Value = 0.2861925343043439
$ws.Cells.Item(10, 16).Value = 0.2861925343043439
$ws.Cells.Item(10, 17).Value = 9.738602310209336
$ws.Cells.Item(10, 18).Value = 87.64742079188402
$ws.Cells.Item(10, 19).Value = 0.01353156587780629
$ws.Cells.Item(10, 20).Value = 0.01353156587780628
$ws.Cells.Item(11, 7).Value = 1.696588
$ws.Cells.Item(11, 8).Value = 5.089764000000001
$ws.Cells.Item(11, 9).Value = 0.04728133775640876
$ws.Cells.Item(11, 10).Value = 0.04728133775640876
$ws.Cells.Item(11, 15).Value = 0.2917347240316885
$ws.Cells.Item(11, 16).Value = 0.2917347240316885
$ws.Cells.Item(11, 17).Value = 9.927192770172001
$ws.Cells.Item(11, 18).Value = 89.344734931548
$ws.Cells.Item(11, 19).Value = 0.01379360802221497
$ws.Cells.Item(11, 20).Value = 0.01379360802221497
$ws.Cells.Item(12, 7).Value = 1.696588
$ws.Cells.Item(12, 8).Value = 5.089764000000001
$ws.Cells.Item(12, 9).Value = 0.04728133775640876
$ws.Cells.Item(12, 10).Value = 0.04728133775640876
$ws.Cells.Item(12, 13).Value = 6.759986
$ws.Cells.Item(12, 14).Value = 20.279958
$ws.Cells.Item(12, 15).Value = 0.3370418707750538
$ws.Cells.Item(12, 16).Value = 0.3370418707750538
$ws.Cells.Item(12, 17).Value = 11.468911127768
$ws.Cells.Item(12, 18).Value = 103.220200149912
$ws.Cells.Item(12, 19).Value = 0.0159357905301672
$ws.Cells.Item(12, 20).Value = 0.0159357905301672
$ws.Cells.Item(13, 7).Value = 1.696588
$ws.Cells.Item(13, 8).Value = 5.089764000000001
$ws.Cells.Item(13, 9).Value = 0.04728133775640876
$ws.Cells.Item(13, 10).Value = 0.04728133775640876
$ws.Cells.Item(13, 13).Value = 1.705448333333333
$ws.Cells.Item(13, 14).Value = 5.116345
$ws.Cells.Item(13, 15).Value = 0.0850308708889137
$ws.Cells.Item(13, 16).Value = 0.0850308708889137
$ws.Cells.Item(13, 17).Value = 2.893443176953333
$ws.Cells.Item(13, 18).Value = 26.04098859258
$ws.Cells.Item(13, 19).Value = 0.004020373326220314
$ws.Cells.Item(13, 20).Value = 0.004020373326220314
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.02836866666666667
$ws.Cells.Item(14, 8).Value = 0.085106
$ws.Cells.Item(14, 9).Value = 0.0007905917702857979
$ws.Cells.Item(14, 10).Value = 0.0007905917702857978
$ws.Cells.Item(14, 13).Value = 5.740110333333334
$ws.Cells.Item(14, 14).Value = 17.220331
$ws.Cells.Item(14, 15).Value = 0.2861925343043439
$ws.Cells.Item(14, 16).Value = 0.2861925343043439
$ws.Cells.Item(14, 17).Value = 0.1628392766762222
$ws.Cells.Item(14, 18).Value = 1.465553490086
$ws.Cells.Item(14, 19).Value = 0.0002262614623382502
$ws.Cells.Item(14, 20).Value = 0.0002262614623382502
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.02836866666666667
$ws.Cells.Item(15, 8).Value = 0.085106
$ws.Cells.Item(15, 9).Value = 0.0007905917702857979
$ws.Cells.Item(15, 10).Value = 0.0007905917702857978
$ws.Cells.Item(15, 15).Value = 0.2917347240316885
$ws.Cells.Item(15, 16).Value = 0.2917347240316885
$ws.Cells.Item(15, 17).Value = 0.165992699838
$ws.Cells.Item(15, 18).Value = 1.493934298542
$ws.Cells.Item(15, 19).Value = 0.0002306430719260513
$ws.Cells.Item(15, 20).Value = 0.0002306430719260513
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.02836866666666667
$ws.Cells.Item(16, 8).Value = 0.085106
$ws.Cells.Item(16, 9).Value = 0.0007905917702857979
$ws.Cells.Item(16, 10).Value = 0.0007905917702857978
$ws.Cells.Item(16, 13).Value = 6.759986
$ws.Cells.Item(16, 14).Value = 20.279958
$ws.Cells.Item(16, 15).Value = 0.3370418707750538
$ws.Cells.Item(16, 16).Value = 0.3370418707750538
$ws.Cells.Item(16, 17).Value = 0.1917717895053334
$ws.Cells.Item(16, 18).Value = 1.725946105548
$ws.Cells.Item(16, 19).Value = 0.0002664625292764869
$ws.Cells.Item(16, 20).Value = 0.0002664625292764869
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.02836866666666667
$ws.Cells.Item(17, 8).Value = 0.085106
$ws.Cells.Item(17, 9).Value = 0.0007905917702857979
$ws.Cells.Item(17, 10).Value = 0.0007905917702857978
$ws.Cells.Item(17, 13).Value = 1.705448333333333
$ws.Cells.Item(17, 14).Value = 5.116345
$ws.Cells.Item(17, 15).Value = 0.0850308708889137
$ws.Cells.Item(17, 16).Value = 0.0850308708889137
$ws.Cells.Item(17, 17).Value = 0.04838129528555556
$ws.Cells.Item(17, 18).Value = 0.43543165757
$ws.Cells.Item(17, 19).Value = 0.0000672247067450094
$ws.Cells.Item(17, 20).Value = 0.00006722470674500939
